$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "43×70=3010" "21×26=546"
Replace-Text "37×34=1258" "65×99=6435"
Replace-Text "53×19=1007" "30×26=780"
Replace-Text "81×79=6399" "21×93=1953"
Replace-Text "41×81=3321" "73×80=5840"
Replace-Text "33×20=660" "34×22=748"
Replace-Text "70×41=2870" "95×34=3230"
Replace-Text "56×83=4648" "48×59=2832"
Replace-Text "97×23=2231" "66×42=2772"
Replace-Text "20×22=440" "22×43=946"
Replace-Text "23×93=2139" "39×99=3861"
Replace-Text "90×74=6660" "92×33=3036"
Replace-Text "16×95=1520" "95×47=4465"
Replace-Text "91×87=7917" "28×63=1764"
Replace-Text "77×32=2464" "64×45=2880"
Replace-Text "18×66=1188" "65×39=2535"
Replace-Text "27×21=567" "78×39=3042"
Replace-Text "92×68=6256" "90×81=7290"
Replace-Text "24×40=960" "87×40=3480"
Replace-Text "17×85=1445" "84×13=1092"
Replace-Text "16×56=896" "13×42=546"
Replace-Text "99×82=8118" "93×52=4836"
Replace-Text "85×34=2890" "20×69=1380"
Replace-Text "24×37=888" "90×55=4950"
Replace-Text "44×16=704" "54×87=4698"
